$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Simple value updates in the first four rows ---
$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"
$t.Cell(4, 1).Range.Text = "301"

# --- Insert three new rows after row 4 (before current row 5) ---
$newRow1 = $t.Rows.Add($t.Rows.Item(5))
$newRow1.Cells.Item(1).Range.Text = "0.00001"

$newRow2 = $t.Rows.Add($t.Rows.Item(6))
$newRow2.Cells.Item(1).Range.Text = "0.00053"

$newRow3 = $t.Rows.Add($t.Rows.Item(7))
$newRow3.Cells.Item(1).Range.Text = "0.00016"

# After the inserts, the old rows shifted down by 3:
#   old row5 (0.00004) -> now row8   (unchanged value)
#   old row6 (0.00029) -> now row9   -> 0.00026
#   old row7 (0.00013) -> now row10  -> 0.00034
#   old row8 (0.00003) -> now row11  -> 0.00041
#   old row9 (0.00016) -> now row12  -> 0.05632
#   old row10 (0.00018) -> now row13 -> DELETE
#   old row11 (0.00022) -> now row13 (after prev delete) -> DELETE
#   old row12 (0.01890) -> now row13 (after prev deletes) -> DELETE

$t.Cell(9, 1).Range.Text = "0.00026"
$t.Cell(10, 1).Range.Text = "0.00034"
$t.Cell(11, 1).Range.Text = "0.00041"
$t.Cell(12, 1).Range.Text = "0.05632"

$t.Rows.Item(13).Delete()
$t.Rows.Item(13).Delete()
$t.Rows.Item(13).Delete()

# --- Collapse the final three multi-column rows into single values ---
$lastCount = $t.Rows.Count
$t.Cell($lastCount - 2, 1).Range.Text = "99.92"
$t.Cell($lastCount - 1, 1).Range.Text = "0.06"
$t.Cell($lastCount, 1).Range.Text = "74"
